$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D22").Value = "'98.36"
$ws.Range("D46").Value = "'9.68"

$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("E5").Value = "  +1.20%  "
$ws.Range("E6").Value = "  +1.84%  "
$ws.Range("E7").Value = "  +12.11%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +1.63%  "
$ws.Range("E10").Value = "  -3.30%  "
$ws.Range("E11").Value = "  -6.29%  "
$ws.Range("E12").Value = "  +10.35%  "
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("E15").Value = "  +0.58%  "
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("E17").Value = "  -2.91%  "
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("E19").Value = "  -3.13%  "
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("E21").Value = "  -3.40%  "
$ws.Range("E22").Value = "  +7.07%  "
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("E24").Value = "  +3.38%  "
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("E26").Value = "  -9.68%  "
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E28").Value = "  +1.60%  "
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("E30").Value = "  +26.20%  "
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("E33").Value = "  -1.86%  "
$ws.Range("E34").Value = "  +4.22%  "
$ws.Range("E35").Value = "  -3.92%  "
$ws.Range("E36").Value = "  +5.45%  "
$ws.Range("E37").Value = "  -3.21%  "
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("E39").Value = "  -9.47%  "
$ws.Range("E40").Value = "  -3.79%  "
$ws.Range("E41").Value = "  +9.05%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("E44").Value = "  -1.16%  "
$ws.Range("E45").Value = "  +4.75%  "
$ws.Range("E46").Value = "  +3.92%  "
$ws.Range("E47").Value = "  -9.89%  "
$ws.Range("E48").Value = "  -3.55%  "
$ws.Range("E49").Value = "  -11.28%  "
$ws.Range("E50").Value = "  -2.55%  "
$ws.Range("E51").Value = "  +0.21%  "
